$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 40, shifting all existing rows
# (old rows 40-65) down to rows 42-67.
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(40).Insert()

# New row 40 (weekly update - newest record)
$ws.Range("A40").Value = 10
$ws.Range("B40").Value = "Vega Modelo de Temuco"
$ws.Range("C40").Value = "La Araucanía"
$ws.Range("D40").Value = 44897
$ws.Range("E40").Value = 9
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100103
$ws.Range("H40").Value = "Frutos de hueso (carozo)"
$ws.Range("I40").Value = 100103003
$ws.Range("J40").Value = "Damasco"
$ws.Range("K40").Value = "Castle Brite"
$ws.Range("L40").Value = "Especial"
$ws.Range("M40").Value = 95
$ws.Range("N40").Value = 18000
$ws.Range("O40").Value = 18000
$ws.Range("P40").Value = 18000
$ws.Range("Q40").Value = "$/bandeja 7 kilos"
$ws.Range("R40").Value = "Provincia de Limarí"
$ws.Range("S40").Value = 2571
$ws.Range("T40").Value = 7

# New row 41 (weekly update - second newest record)
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44897
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100103
$ws.Range("H41").Value = "Frutos de hueso (carozo)"
$ws.Range("I41").Value = 100103003
$ws.Range("J41").Value = "Damasco"
$ws.Range("K41").Value = "Castle Brite"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 45
$ws.Range("N41").Value = 28000
$ws.Range("O41").Value = 28000
$ws.Range("P41").Value = 28000
$ws.Range("Q41").Value = "$/bandeja 18 kilos"
$ws.Range("R41").Value = "Provincia de Limarí"
$ws.Range("S41").Value = 1556
$ws.Range("T41").Value = 18
